$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Consolidate "A slide" title runs into a single run.
# Setting the text to a transient, non-overlapping value first forces the
# host's diff-based text-range writer to rebuild the paragraph from a
# single fresh run (instead of reusing/splitting the 3 existing runs), so
# the final assignment collapses "A" + " " + "slide" into one run.
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "-"
$titleRange.Text = "A slide"

# Consolidate "Followed by a picture" caption runs into a single run, same technique.
$captionRange = $s.Shapes.Item(4).TextFrame.TextRange
$captionRange.Text = "-"
$captionRange.Text = "Followed by a picture"
